# latest update Ash Part
# Adds SavAcc/CheckAcc summary columns to Clients, inserts a new ID/Acc#
# label column on Checkin/Saving/Money Market and fills in account-balance
# detail strings on Saving and Money Market.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Clients sheet: new SavAcc / CheckAcc columns (F:G) + a PIN correction
# ---------------------------------------------------------------------
$clients = $wb.Worksheets.Item("Clients")

$clients.Range("F1").Value2 = "SavAcc"
$clients.Range("G1").Value2 = "CheckAcc"

$clients.Range("B3").Value2 = "0"

$clients.Range("F2").Value2 = "1"
$clients.Range("G2").Value2 = "1"

$clients.Range("F3").Value2 = "3"
$clients.Range("G3").Value2 = "3"

$clients.Range("F4").Value2 = "1"
$clients.Range("G4").Value2 = "1"

$clients.Range("F5").Value2 = "0"
$clients.Range("G5").Value2 = "0"

$clients.Range("F6").Value2 = 0
$clients.Range("G6").Value2 = 0

$clients.Range("F7").Value2 = 0
$clients.Range("G7").Value2 = 0

$clients.Range("F8").Value2 = 0
$clients.Range("G8").Value2 = 0

$clients.Range("F9").Value2 = 0
$clients.Range("G9").Value2 = 0

$clients.Range("F10").Value2 = 0
$clients.Range("G10").Value2 = 0

$clients.Range("F11").Value2 = 0
$clients.Range("G11").Value2 = 0

# ---------------------------------------------------------------------
# Checkin / Saving / Money Market: insert a leading ID / Acc#n label
# column, pushing the existing account data one column to the right.
# ---------------------------------------------------------------------
$checkin = $wb.Worksheets.Item("Checkin")
$saving = $wb.Worksheets.Item("Saving")
$money = $wb.Worksheets.Item("Money Market")

$checkin.Columns.Item(1).Insert()
$saving.Columns.Item(1).Insert()
$money.Columns.Item(1).Insert()

foreach ($ws in @($checkin, $saving, $money)) {
    $ws.Range("A1").Value2 = "ID"
    $ws.Range("A2").Value2 = "Acc#1"
    $ws.Range("A3").Value2 = "Acc#2"
    $ws.Range("A4").Value2 = "Acc#3"
    $ws.Range("A5").Value2 = "Acc#4"
    $ws.Range("A6").Value2 = "Acc#5"
}

# ---------------------------------------------------------------------
# Saving sheet: account/balance detail strings
# ---------------------------------------------------------------------
$saving.Range("D2").Value2 = "204203#3"
$saving.Range("E2").Value2 = "204204#3"
$saving.Range("F2").Value2 = "204205#3"
$saving.Range("G2").Value2 = "204206#3"
$saving.Range("H2").Value2 = "204207#3"
$saving.Range("I2").Value2 = "204208#3"
$saving.Range("J2").Value2 = "204209#3"
$saving.Range("K2").Value2 = "204210#3"
$saving.Range("F3").Value2 = "204212#1"
$saving.Range("I3").Value2 = "204213#3"
$saving.Range("F4").Value2 = "204214#3"
$saving.Range("B2").Value2 = "2042013874#3"
$saving.Range("C2").Value2 = "204202#0#3"
$saving.Range("B3").Value2 = "204211#23#2"

# ---------------------------------------------------------------------
# Money Market sheet: account/balance detail strings
# ---------------------------------------------------------------------
$money.Range("B2").Value2 = "304201#543"
$money.Range("C2").Value2 = "304202#123456"
$money.Range("D2").Value2 = "304203#4354354"
$money.Range("E2").Value2 = "304204#40000"
$money.Range("F2").Value2 = "304205#2300"
$money.Range("G2").Value2 = "304206#0"
$money.Range("H2").Value2 = "304207#12"
$money.Range("I2").Value2 = "304208#45000"
$money.Range("J2").Value2 = "304209#2000000"
$money.Range("K2").Value2 = "304210#323232"

# ---------------------------------------------------------------------
# Selection / active-sheet bookkeeping to mirror the saved view state.
# ---------------------------------------------------------------------
$checkin.Range("A1:A1048576").Select()
$saving.Range("A1:A1048576").Select()
$money.Range("C23").Select()
$clients.Range("G15").Select()
$clients.Select()
